$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2: was a plain value (3); now driven by a formula 2*1.5 (still 3)
$ws.Range("E2").Formula = "=2*1.5"

# E4: was a plain value (7.5); now driven by a formula 4*1.5 (=> 6)
$ws.Range("E4").Formula = "=4*1.5"

# Row 5: fill in Ranking / Points for ranking / Points for delta / Points for bonus
$ws.Range("B5").Value = 40
$ws.Range("C5").Formula = "=(6.7+6.9)/2"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0

# Row 6: fill in Ranking / Points for ranking / Points for delta / Points for bonus
$ws.Range("B6").Formula = "=40"
$ws.Range("C6").Formula = "=6.8"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0

# Row 7: fill in Ranking / Points for ranking / Points for delta / Points for bonus
$ws.Range("B7").Value = 33
$ws.Range("C7").Formula = "=17-16*LOG((B7-1)/(70-1) + 1,2)"
$ws.Range("D7").Formula = "=1"
$ws.Range("E7").Formula = "=1.5*1"

# F7's average formula becomes a weighted SUM that double counts the latest ranking
$ws.Range("F7").Formula = "=SUM(C2:C7,C7)/7"

# Move the active selection to I7
$ws.Range("I7").Select() | Out-Null
